# Part 6 README edit: expand the "Trello" bullet in the speaker notes of
# slide 2 ("Tools") with the additional commentary about why the tracker
# wasn't more useful to the team.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$notes = $s.NotesPage
$shape = $notes.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$para1 = "Trello: it did what it wanted but it wasn" + [char]8217 + "t exactly what we would have wanted. Maybe it was just the way we used it but the stories needed to be smaller and we didn" + [char]8217 + "t actually use it because it was more convenient for us to talk or group text about what we had to do, doing, done, etc. Maybe if we had a bigger group or were part of a company the tracker would have been more useful for us"
$para2 = "Git: would have been really useful to have notifications"

$tr.Text = $para1 + [char]10 + $para2
